$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corregir nombres (apellidos) de los alumnos
$ws.Range("A4").Value = "Acziendo"
$ws.Range("A3").Value = "Torreta"
$ws.Range("A5").Value = "Finolla"
$ws.Range("A6").Value = "Serf"
$ws.Range("A7").Value = "Pirez"
$ws.Range("A8").Value = "Luciendo"

# Subrayar el apellido final (Martinez) en A9
$ws.Range("A9").Font.Underline = $true

# Configuracion de pagina (tamano carta/A4 y orientacion vertical)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selecciona A9 como celda activa final
[void]$ws.Range("A9").Select()
